{"js": "// Update the PA1 due date from \"Friday, May 11, 2018\" to\n// \"Wednesday, February 6, 2019\" (commit: \"updated pa due date (thanks #2)\").\nconst body = context.document.body;\n\nconst oldDate = \"Friday, May 11, 2018\";\nconst newDate = \"Wednesday, February 6, 2019\";\n\nconst results = body.search(oldDate, {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newDate, \"Replace\");\n}\nawait context.sync();\n\n// Word tracks the location of the most recent edit with the hidden\n// \"_GoBack\" bookmark. Move it to just after the newly typed date (where\n// the user's cursor would have been left), matching what Word itself does\n// when you type new text into a document. A document can only have one\n// bookmark with a given name, so drop the old one first.\nconst newDateRanges = body.search(newDate, { matchCase: true });\nnewDateRanges.load(\"items\");\nawait context.sync();\n\nif (newDateRanges.items.length > 0) {\n  context.document.deleteBookmark(\"_GoBack\");\n  const afterDate = newDateRanges.items[0].getRange(\"End\");\n  afterDate.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Update the PA1 due date from \"Friday, May 11, 2018\" to\n# \"Wednesday, February 6, 2019\" (commit: \"updated pa due date (thanks #2)\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Friday, May 11, 2018\"\n$find.Replacement.Text = \"Wednesday, February 6, 2019\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Word tracks the location of the most recent edit with the hidden\n# \"_GoBack\" bookmark. Move it to just after the newly typed date (where\n# the user's cursor would have been left), matching what Word itself\n# does when you type new text into a document.\n$goBack = $d.Content.Find\n$goBack.ClearFormatting()\n$goBack.Text = \"Wednesday, February 6, 2019\"\n$goBack.Forward = $true\n$goBack.Wrap = 1\nif ($goBack.Execute()) {\n    $rng = $goBack.Parent.Duplicate\n    $rng.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n}\n"}
